# Add data for 2022-07-02:
# - Sheet tab renamed from "Through 2022-06-23" to "Through 2022-06-24"
# - Header cell B1 text updated to match ("June 2022 (through June 23)" -> "...June 24)")
# - A handful of cells in the "current month" (and a few other) columns are
#   incremented/added to reflect the newly-ingested carjacking records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "Through 2022-06-24"

# Update the column header text (shared string) to match the new "through" date.
$ws.Range("B1").Value = "June 2022 (through June 24)"

# Updated / incremented counts.
$ws.Range("T2").Value = 4
$ws.Range("AL2").Value = 3

$ws.Range("Z3").Value = 2
$ws.Range("AF3").Value = 2
$ws.Range("AR3").Value = 2

$ws.Range("H4").Value = 6

$ws.Range("AF5").Value = 8

$ws.Range("AR6").Value = 1

$ws.Range("N9").Value = 6

$ws.Range("AF10").Value = 4

$ws.Range("B11").Value = 2

$ws.Range("H14").Value = 12
$ws.Range("AL14").Value = 2

$ws.Range("H17").Value = 3

$ws.Range("N27").Value = 2

$ws.Range("B28").Value = 1

$ws.Range("AL37").Value = 1

$ws.Range("Z65").Value = 3

$ws.Range("B71").Value = 3

$ws.Range("B81").Value = 2

$ws.Range("H85").Value = 3

$ws.Range("B91").Value = 1

$ws.Range("AL92").Value = 1
